$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds order numbers that look numeric ("11149", "11201", ...).
# The source file stores them as literal text, so force text formatting on
# that column before writing the values - otherwise Excel's type-sniffing
# on `.Value` would silently coerce them to numbers. Column B's dates
# ("19/11/2024", "21/11/2024", ...) and columns C-E are never numeric-
# looking, so they don't need this treatment.
$ws.Range("A2:A18").NumberFormat = "@"

$data = @{
  2  = @("11149","19/11/2024","manhã","Desconhecido","Distrito Industrial")
  3  = @("11201","19/11/2024","manhã","Desconhecido","Ponte do Imaruim")
  4  = @("11142","19/11/2024","manhã","Desconhecido","Centro")
  5  = @("11179","19/11/2024","manhã","Desconhecido","Nova Palhoça")
  6  = @("11064","19/11/2024","manhã","Desconhecido","Barra do Aririú")
  7  = @("11158","19/11/2024","tarde","Desconhecido","Forquilhinha")
  8  = @("11151","19/11/2024","tarde","Desconhecido","Forquilhas")
  9  = @("11222","19/11/2024","tarde","Desconhecido","Sertão do Maruim")
  10 = @("11194","19/11/2024","tarde","Desconhecido","Pedra Branca")
  11 = @("11193","19/11/2024","tarde","Desconhecido","Pedra Branca")
  12 = @("11197","21/11/2024","manhã","Desconhecido","Picadas do Sul")
  13 = @("11140","21/11/2024","manhã","Desconhecido","Centro")
  14 = @("11071","21/11/2024","manhã","Desconhecido","Barreiros")
  15 = @("11177","21/11/2024","manhã","Desconhecido","Nossa Senhora do Rosário")
  16 = @("11200","21/11/2024","tarde","Desconhecido","Ponta de Baixo")
  17 = @("11141","21/11/2024","tarde","Desconhecido","Centro")
  18 = @("11144","21/11/2024","tarde","Desconhecido","Centro")
}

foreach ($r in $data.Keys) {
  $row = $data[$r]
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
}

# The source data was trimmed from 20 rows to 18 - remove the two now-
# obsolete trailing rows (old 11088/Cachoeiras and 11105/Campeche entries).
$ws.Range("A19:E20").Delete()

# Column E (Bairro) got a bit narrower in the new layout (stored width 29 -> 26).
# The engine quantizes ColumnWidth to 1/6-character steps on write/round-trip,
# so feed it a value from the middle of the input band that resolves to 26
# after quantization, rather than 26 itself (which resolves to ~26.83).
$ws.Columns.Item(5).ColumnWidth = 25.1667
